# Auto-generated script to rewrite sheet1 data grid to match target state
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing used range content first (rows 1-17, columns A-Z)
$ws.Range("A1:Z17").ClearContents()

# Row 1
$ws.Range("A1").Value = 'ReachName'
$ws.Range("B1").Value = 'Basin'
$ws.Range("C1").Value = 'Assessment.Unit'
$ws.Range("D1").Value = 'Spring.Chinook.Reach'
$ws.Range("E1").Value = 'Steelhead.Reach'
$ws.Range("F1").Value = 'Bull.Trout.Reach'
$ws.Range("G1").Value = 'BankStability_score'
$ws.Range("H1").Value = 'ChannelStability_score'
$ws.Range("I1").Value = 'Stability_Mean'
$ws.Range("J1").Value = 'CoarseSubstrate_score'
$ws.Range("K1").Value = 'Cover-Wood_score'
$ws.Range("L1").Value = 'Flow-SummerBaseFlow_score'
$ws.Range("M1").Value = 'FloodplainConnectivity_score'
$ws.Range("N1").Value = 'Off-Channel/Side-Channels_score'
$ws.Range("O1").Value = 'PoolQuantity&Quality_score'
$ws.Range("P1").Value = 'Riparian-CanopyCover_score'
$ws.Range("Q1").Value = 'Riparian-Disturbance_score'
$ws.Range("R1").Value = 'Riparian_Mean'
$ws.Range("S1").Value = 'Temperature-Rearing_score'
$ws.Range("T1").Value = 'HQ_Sum'
$ws.Range("U1").Value = 'HQ_Pct'
$ws.Range("V1").Value = 'HQ_Score_Restoration'
$ws.Range("W1").Value = 'HQ_Score_Protection'
$ws.Range("X1").Value = 'unacceptable_1_indiv_habitat_attributes'
$ws.Range("Y1").Value = 'at_risk_2_or_3_indiv_habitat_attributes'
$ws.Range("Z1").Value = 'unacceptable_AND_at_risk_1_to_3_indiv_habitat_attributes'

# Row 2
$ws.Range("A2").Value = 'Entiat River Lake 02'
$ws.Range("B2").Value = 'Entiat'
$ws.Range("C2").Value = 'Entiat River-Lake Creek'
$ws.Range("D2").Value = 'yes'
$ws.Range("E2").Value = 'yes'
$ws.Range("F2").Value = 'yes'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 3
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 5
$ws.Range("T2").Value = 27
$ws.Range("U2").Value = 0.6
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 'Cover-Wood,FloodplainConnectivity'
$ws.Range("Y2").Value = 'Stability,Flow-SummerBaseFlow,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z2").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'

# Row 3
$ws.Range("A3").Value = 'Entiat River Lake 04'
$ws.Range("B3").Value = 'Entiat'
$ws.Range("C3").Value = 'Entiat River-Lake Creek'
$ws.Range("D3").Value = 'yes'
$ws.Range("E3").Value = 'yes'
$ws.Range("F3").Value = 'yes'
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 5
$ws.Range("R3").Value = 4
$ws.Range("S3").Value = 5
$ws.Range("T3").Value = 32
$ws.Range("U3").Value = 0.7111111111111111
$ws.Range("V3").Value = 5
$ws.Range("W3").Value = 3
$ws.Range("X3").Value = 'Cover-Wood,PoolQuantity&Quality'
$ws.Range("Y3").Value = 'Flow-SummerBaseFlow,Off-Channel/Side-Channels'
$ws.Range("Z3").Value = 'Cover-Wood,Flow-SummerBaseFlow,Off-Channel/Side-Channels,PoolQuantity&Quality'

# Row 4
$ws.Range("A4").Value = 'Entiat River Potato 07'
$ws.Range("B4").Value = 'Entiat'
$ws.Range("C4").Value = 'Entiat River-Potato Creek'
$ws.Range("D4").Value = 'yes'
$ws.Range("E4").Value = 'yes'
$ws.Range("F4").Value = 'yes'
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 3
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 31
$ws.Range("U4").Value = 0.6888888888888889
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 1
$ws.Range("Y4").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian,Temperature-Rearing'
$ws.Range("Z4").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian,Temperature-Rearing'

# Row 5
$ws.Range("A5").Value = 'Entiat River Potato 08'
$ws.Range("B5").Value = 'Entiat'
$ws.Range("C5").Value = 'Entiat River-Potato Creek'
$ws.Range("D5").Value = 'yes'
$ws.Range("E5").Value = 'yes'
$ws.Range("F5").Value = 'yes'
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 5
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 3
$ws.Range("N5").Value = 3
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 33
$ws.Range("U5").Value = 0.7333333333333333
$ws.Range("V5").Value = 5
$ws.Range("W5").Value = 3
$ws.Range("Y5").Value = 'Stability,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian,Temperature-Rearing'
$ws.Range("Z5").Value = 'Stability,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian,Temperature-Rearing'

# Row 6
$ws.Range("A6").Value = 'Nason Creek Lower 03'
$ws.Range("B6").Value = 'Wenatchee'
$ws.Range("C6").Value = 'Lower Nason Creek'
$ws.Range("D6").Value = 'yes'
$ws.Range("E6").Value = 'yes'
$ws.Range("F6").Value = 'yes'
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = 5
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 5
$ws.Range("R6").Value = 4
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 35
$ws.Range("U6").Value = 0.7777777777777778
$ws.Range("V6").Value = 5
$ws.Range("W6").Value = 3
$ws.Range("X6").Value = 'Temperature-Rearing'
$ws.Range("Y6").Value = 'Flow-SummerBaseFlow,FloodplainConnectivity'
$ws.Range("Z6").Value = 'Flow-SummerBaseFlow,FloodplainConnectivity,Temperature-Rearing'

# Row 7
$ws.Range("A7").Value = 'Nason Creek Lower 04'
$ws.Range("B7").Value = 'Wenatchee'
$ws.Range("C7").Value = 'Lower Nason Creek'
$ws.Range("D7").Value = 'yes'
$ws.Range("E7").Value = 'yes'
$ws.Range("F7").Value = 'yes'
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 5
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = 3
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 3
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 27
$ws.Range("U7").Value = 0.6
$ws.Range("V7").Value = 5
$ws.Range("W7").Value = 1
$ws.Range("X7").Value = 'PoolQuantity&Quality,Temperature-Rearing'
$ws.Range("Y7").Value = 'Stability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z7").Value = 'Stability,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 8
$ws.Range("A8").Value = 'Nason Creek Lower 05'
$ws.Range("B8").Value = 'Wenatchee'
$ws.Range("C8").Value = 'Lower Nason Creek'
$ws.Range("D8").Value = 'yes'
$ws.Range("E8").Value = 'yes'
$ws.Range("F8").Value = 'yes'
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 5
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 3
$ws.Range("N8").Value = 3
$ws.Range("O8").Value = 1
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 2
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 24
$ws.Range("U8").Value = 0.5333333333333333
$ws.Range("V8").Value = 5
$ws.Range("W8").Value = 1
$ws.Range("X8").Value = 'PoolQuantity&Quality,Temperature-Rearing'
$ws.Range("Y8").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z8").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 9
$ws.Range("A9").Value = 'Nason Creek Lower 06'
$ws.Range("B9").Value = 'Wenatchee'
$ws.Range("C9").Value = 'Lower Nason Creek'
$ws.Range("D9").Value = 'yes'
$ws.Range("E9").Value = 'yes'
$ws.Range("F9").Value = 'yes'
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 3
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 3
$ws.Range("S9").Value = 1
$ws.Range("T9").Value = 25
$ws.Range("U9").Value = 0.5555555555555556
$ws.Range("V9").Value = 5
$ws.Range("W9").Value = 1
$ws.Range("X9").Value = 'PoolQuantity&Quality,Temperature-Rearing'
$ws.Range("Y9").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z9").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 10
$ws.Range("A10").Value = 'Nason Creek Lower 07'
$ws.Range("B10").Value = 'Wenatchee'
$ws.Range("C10").Value = 'Lower Nason Creek'
$ws.Range("D10").Value = 'yes'
$ws.Range("E10").Value = 'yes'
$ws.Range("F10").Value = 'yes'
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 3
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 3
$ws.Range("M10").Value = 3
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 1
$ws.Range("P10").Value = 3
$ws.Range("Q10").Value = 3
$ws.Range("R10").Value = 3
$ws.Range("S10").Value = 1
$ws.Range("T10").Value = 25
$ws.Range("U10").Value = 0.5555555555555556
$ws.Range("V10").Value = 5
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 'PoolQuantity&Quality,Temperature-Rearing'
$ws.Range("Y10").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z10").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 11
$ws.Range("A11").Value = 'Nason Creek Lower 09'
$ws.Range("B11").Value = 'Wenatchee'
$ws.Range("C11").Value = 'Lower Nason Creek'
$ws.Range("D11").Value = 'yes'
$ws.Range("E11").Value = 'yes'
$ws.Range("F11").Value = 'yes'
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 5
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 1
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 3
$ws.Range("P11").Value = 1
$ws.Range("Q11").Value = 3
$ws.Range("R11").Value = 2
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 19
$ws.Range("U11").Value = 0.4222222222222222
$ws.Range("V11").Value = 5
$ws.Range("W11").Value = 1
$ws.Range("X11").Value = 'Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Temperature-Rearing'
$ws.Range("Y11").Value = 'Stability,Flow-SummerBaseFlow,PoolQuantity&Quality,Riparian'
$ws.Range("Z11").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 12
$ws.Range("A12").Value = 'Nason Creek Lower 10'
$ws.Range("B12").Value = 'Wenatchee'
$ws.Range("C12").Value = 'Lower Nason Creek'
$ws.Range("D12").Value = 'yes'
$ws.Range("E12").Value = 'yes'
$ws.Range("F12").Value = 'yes'
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 5
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 3
$ws.Range("M12").Value = 1
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 3
$ws.Range("P12").Value = 1
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 1
$ws.Range("T12").Value = 19
$ws.Range("U12").Value = 0.4222222222222222
$ws.Range("V12").Value = 5
$ws.Range("W12").Value = 1
$ws.Range("X12").Value = 'Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,Temperature-Rearing'
$ws.Range("Y12").Value = 'Stability,Flow-SummerBaseFlow,PoolQuantity&Quality,Riparian'
$ws.Range("Z12").Value = 'Stability,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 13
$ws.Range("A13").Value = 'Nason Creek Lower 11'
$ws.Range("B13").Value = 'Wenatchee'
$ws.Range("C13").Value = 'Lower Nason Creek'
$ws.Range("D13").Value = 'yes'
$ws.Range("E13").Value = 'yes'
$ws.Range("F13").Value = 'yes'
$ws.Range("G13").Value = 3
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = 3
$ws.Range("N13").Value = 3
$ws.Range("O13").Value = 1
$ws.Range("P13").Value = 3
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 2
$ws.Range("S13").Value = 1
$ws.Range("T13").Value = 22
$ws.Range("U13").Value = 0.4888888888888889
$ws.Range("V13").Value = 5
$ws.Range("W13").Value = 1
$ws.Range("X13").Value = 'PoolQuantity&Quality,Temperature-Rearing'
$ws.Range("Y13").Value = 'Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,Riparian'
$ws.Range("Z13").Value = 'Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 14
$ws.Range("A14").Value = 'Nason Creek Lower 12'
$ws.Range("B14").Value = 'Wenatchee'
$ws.Range("C14").Value = 'Lower Nason Creek'
$ws.Range("D14").Value = 'yes'
$ws.Range("E14").Value = 'yes'
$ws.Range("F14").Value = 'yes'
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 3
$ws.Range("M14").Value = 3
$ws.Range("N14").Value = 3
$ws.Range("O14").Value = 3
$ws.Range("P14").Value = 1
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 21
$ws.Range("U14").Value = 0.4666666666666667
$ws.Range("V14").Value = 5
$ws.Range("W14").Value = 1
$ws.Range("X14").Value = 'Cover-Wood,Riparian,Temperature-Rearing'
$ws.Range("Y14").Value = 'Stability,CoarseSubstrate,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality'
$ws.Range("Z14").Value = 'Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 15
$ws.Range("A15").Value = 'Nason Creek Lower 13'
$ws.Range("B15").Value = 'Wenatchee'
$ws.Range("C15").Value = 'Lower Nason Creek'
$ws.Range("D15").Value = 'yes'
$ws.Range("E15").Value = 'yes'
$ws.Range("F15").Value = 'yes'
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 3
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 3
$ws.Range("M15").Value = 3
$ws.Range("N15").Value = 3
$ws.Range("O15").Value = 3
$ws.Range("P15").Value = 1
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = 1
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 21
$ws.Range("U15").Value = 0.4666666666666667
$ws.Range("V15").Value = 5
$ws.Range("W15").Value = 1
$ws.Range("X15").Value = 'Cover-Wood,Riparian,Temperature-Rearing'
$ws.Range("Y15").Value = 'Stability,CoarseSubstrate,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality'
$ws.Range("Z15").Value = 'Stability,CoarseSubstrate,Cover-Wood,Flow-SummerBaseFlow,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Row 16
$ws.Range("A16").Value = 'Nason Creek Lower 15'
$ws.Range("B16").Value = 'Wenatchee'
$ws.Range("C16").Value = 'Lower Nason Creek'
$ws.Range("D16").Value = 'yes'
$ws.Range("E16").Value = 'yes'
$ws.Range("F16").Value = 'yes'
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 3
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 3
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 3
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 3
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 27
$ws.Range("U16").Value = 0.6
$ws.Range("V16").Value = 5
$ws.Range("W16").Value = 1
$ws.Range("X16").Value = 'Temperature-Rearing'
$ws.Range("Y16").Value = 'Stability,CoarseSubstrate,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian'
$ws.Range("Z16").Value = 'Stability,CoarseSubstrate,Cover-Wood,FloodplainConnectivity,Off-Channel/Side-Channels,PoolQuantity&Quality,Riparian,Temperature-Rearing'

# Remove now-stale row 17 so used range/dimension becomes A1:Z16
$ws.Rows.Item(17).Delete()
